$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing header style (bold, centered, bordered) from A1
# before overwriting values, then reapply it to the full A1:AG1 header range.
$ws.Cells.Item(1, 1).Copy()

# Set the final header row (new columns inserted, old ones reordered/renamed).
$ws.Range("A1").Value = "Arena Id"
$ws.Range("B1").Value = "Backboard Shakes"
$ws.Range("C1").Value = "City Floor Id"
$ws.Range("D1").Value = "City Short Name"
$ws.Range("E1").Value = "Crowd Type"
$ws.Range("F1").Value = "DornasID"
$ws.Range("G1").Value = "Floor Id"
$ws.Range("H1").Value = "ID"
$ws.Range("I1").Value = "Stadium Vitals - LED ID"
$ws.Range("J1").Value = "Stadium Vitals - ALT_FLOOR_FILE#1"
$ws.Range("K1").Value = "Stadium Vitals - ALT_FLOOR_FILE#2"
$ws.Range("L1").Value = "Stadium Vitals - ALT_FLOOR_FILE#3"
$ws.Range("M1").Value = "Stadium Vitals - ALT_FLOOR_JERSEY#1"
$ws.Range("N1").Value = "Stadium Vitals - ALT_FLOOR_JERSEY#2"
$ws.Range("O1").Value = "Stadium Vitals - ALT_FLOOR_JERSEY#3"
$ws.Range("P1").Value = "Stadium Vitals - ARENA_CAPACITY"
$ws.Range("Q1").Value = "Stadium Vitals - ARENA_FILE"
$ws.Range("R1").Value = "Stadium Vitals - BASKETBALL_TYPE"
$ws.Range("S1").Value = "Stadium Vitals - CITY_ABB"
$ws.Range("T1").Value = "Stadium Vitals - CITY_NAME"
$ws.Range("U1").Value = "Stadium Vitals - CROWD_LOUDNESS"
$ws.Range("V1").Value = "Stadium Vitals - DORNAS_ID"
$ws.Range("W1").Value = "Stadium Vitals - FLOOR_FILE"
$ws.Range("X1").Value = "Stadium Vitals - HOME_BASKET"
$ws.Range("Y1").Value = "Stadium Vitals - LOCATION_TYPE"
$ws.Range("Z1").Value = "Stadium Vitals - MUSIC_PROBABILITY"
$ws.Range("AA1").Value = "Stadium Vitals - NAME"
$ws.Range("AB1").Value = "Stadium Vitals - NICKNAME"
$ws.Range("AC1").Value = "Stadium Vitals - SOUND_ID"
$ws.Range("AD1").Value = "Stadium Vitals - STATE"
$ws.Range("AE1").Value = "Stadium Vitals - TYPE"
$ws.Range("AF1").Value = "Stadium Vitals - UNIQUEID"
$ws.Range("AG1").Value = "State Short Name"

# Re-apply the header formatting (font, alignment, borders) across the full range,
# including the newly added columns Z:AG, to match the original header style.
$ws.Range("A1:AG1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
